# Backlog_7.xlsx edit
#
# Semana 07 has passed, so the "Semana" column (C) no longer needs the
# "Semana 07" text label -- it is replaced everywhere by the plain week
# number 7. Once every usage of the "Semana 07" shared string is gone,
# saving naturally drops it from the shared-string table and re-indexes
# every other shared string used by the Responsavel (B) / Status (I)
# columns accordingly.
#
# Also updates the active sheet/selection to match where work left off:
# SPN!F18 is now the live cell (and the active tab), while ITI keeps its
# own last selection (C2:C29) for when it is revisited.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("SPN")
$ws2 = $wb.Worksheets.Item("ITI")

# --- SPN: rows 2-13 ---
for ($r = 2; $r -le 13; $r++) {
    $ws1.Cells.Item($r, 3).Value = 7
}

# --- ITI: rows 2-29 ---
for ($r = 2; $r -le 29; $r++) {
    $ws2.Cells.Item($r, 3).Value = 7
}

# Remember each sheet's own selection.
$ws2.Range("C2:C29").Select() | Out-Null
$ws1.Range("F18").Select() | Out-Null

# SPN is the sheet left on screen (tab 1 / index 0).
$ws1.Activate()
